{"js": "// Change the East Asian / complex-script font fallbacks used by the\n// document's paragraph styles:\n//   - Normal / Heading: eastAsia font DejaVu Sans -> Tahoma\n//   - List / Caption / Index: explicitly pin the complex-script (w:cs)\n//     font to DejaVu Sans (previously inherited / unset)\n//\n// Note: the document's docDefaults (<w:docDefaults>/<w:rPrDefault>) is\n// not reachable through the Word JS API (there is no scriptable object\n// for it, same as in real Word), so it is intentionally left untouched\n// here; only the named styles below are addressable.\n\nconst styles = context.document.getStyles();\n\nconst normal = styles.getByNameOrNullObject(\"Normal\");\nconst heading = styles.getByNameOrNullObject(\"Heading\");\nconst list = styles.getByNameOrNullObject(\"List\");\nconst caption = styles.getByNameOrNullObject(\"Caption\");\nconst index = styles.getByNameOrNullObject(\"Index\");\n\nnormal.load(\"isNullObject\");\nheading.load(\"isNullObject\");\nlist.load(\"isNullObject\");\ncaption.load(\"isNullObject\");\nindex.load(\"isNullObject\");\nawait context.sync();\n\nif (!normal.isNullObject) {\n  normal.font.nameFarEast = \"Tahoma\";\n}\nif (!heading.isNullObject) {\n  heading.font.nameFarEast = \"Tahoma\";\n}\nif (!list.isNullObject) {\n  list.font.nameBidirectional = \"DejaVu Sans\";\n}\nif (!caption.isNullObject) {\n  caption.font.nameBidirectional = \"DejaVu Sans\";\n}\nif (!index.isNullObject) {\n  index.font.nameBidirectional = \"DejaVu Sans\";\n}\n\nawait context.sync();\n", "ps1": "# Change the East Asian / complex-script font fallbacks used by the\n# document's paragraph styles:\n#   - Normal / Heading: eastAsia font DejaVu Sans -> Tahoma\n#   - List / Caption / Index: explicitly pin the complex-script (w:cs)\n#     font to DejaVu Sans (previously inherited / unset)\n#\n# Note: the document's docDefaults (<w:docDefaults>/<w:rPrDefault>) is\n# not reachable through the Word object model (there is no Styles()\n# entry / property for it, same as in real Word automation), so it is\n# intentionally left untouched here; only the named styles below are\n# addressable.\n\n$d = $word.ActiveDocument\n\n$d.Styles(\"Normal\").Font.NameFarEast = \"Tahoma\"\n$d.Styles(\"Heading\").Font.NameFarEast = \"Tahoma\"\n\n$d.Styles(\"List\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles(\"Caption\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles(\"Index\").Font.NameBi = \"DejaVu Sans\"\n"}
